$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 - best_params: update hyperparameter values for a few models
$ws.Range("F2").Value = "{'max_depth': 50, 'n_estimators': 50}"
$ws.Range("G2").Value = "{'learning_rate': 0.1, 'max_depth': 3, 'n_estimators': 100}"
$ws.Range("H2").Value = "{'learning_rate': 1, 'n_estimators': 100}"
$ws.Range("K2").Value = "{'activation': 'leaky_relu', 'b_random_vec_range': [0, 10], 'lam': 2, 'n_layer': 16, 'n_nodes': 128, 'random_seed': 882, 'same_feature': True, 'w_random_vec_range': [-10, 10]}"

# Row 3 - rmse
$ws.Range("B3").Value = 0.1078316630187899
$ws.Range("C3").Value = 0.1200915680019607
$ws.Range("D3").Value = 0.1157733381862321
$ws.Range("E3").Value = 0.08532260496695201
$ws.Range("F3").Value = 0.07994515816766538
$ws.Range("G3").Value = 0.08277524828396099
$ws.Range("H3").Value = 0.09537684150853717
$ws.Range("I3").Value = 0.08381964103540558
$ws.Range("J3").Value = 0.08341334576730942
$ws.Range("K3").Value = 0.03765903314938047

# Row 4 - r2
$ws.Range("B4").Value = 0.2233501807304878
$ws.Range("C4").Value = 0.078317987091215
$ws.Range("D4").Value = 0.1312474590069004
$ws.Range("E4").Value = 0.5162473323301832
$ws.Range("F4").Value = 0.5760342119476762
$ws.Range("G4").Value = 0.5495628755824622
$ws.Range("H4").Value = 0.4025244311519659
$ws.Range("I4").Value = 0.5138205292733218
$ws.Range("J4").Value = 0.5434218072769883
$ws.Range("K4").Value = 0.8832398449700054

# Row 5 - mape
$ws.Range("B5").Value = 48.60246420385555
$ws.Range("C5").Value = 65.30119490085931
$ws.Range("D5").Value = 62.27375119029092
$ws.Range("E5").Value = 39.15986619907513
$ws.Range("F5").Value = 30.93111460380483
$ws.Range("G5").Value = 31.56129552824202
$ws.Range("H5").Value = 44.53880581129712
$ws.Range("I5").Value = 29.48935885561542
$ws.Range("J5").Value = 32.25064957194685
$ws.Range("K5").Value = 11.51816879022205
